$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (cohort 2018, period_index 4): num_customers 50 -> 51, retention_rate recalculated
$ws.Range("C27").Value = 51
$ws.Range("E27").Value = 0.02264653641207815

# Row 31 (cohort 2019, period_index 3): num_customers 53 -> 54, retention_rate recalculated
$ws.Range("C31").Value = 54
$ws.Range("E31").Value = 0.02335640138408305

# Row 37 (cohort 2024, period_index 0): num_customers 917 -> 923, cohort_size 917 -> 923
$ws.Range("C37").Value = 923
$ws.Range("D37").Value = 923
